$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old D:K -> E:L),
# matching the new fiscal-year column added to each of the three
# statements (Income Statement, Balance Sheet, Cash Flow).
$ws.Range("D1").EntireColumn.Insert()

# Carry the number formats/styles from the (old) neighboring column,
# now in column E, into the freshly inserted column D, per block of
# rows (skipping the blank separator / section-header rows that have
# no D:K data at all).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate the new column D with the new period's figures.

# --- Income Statement (rows 7-35) ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 560500
$ws.Range("D9").Value = 263700
$ws.Range("D10").Value = 296700
$ws.Range("D12").Value = 19300
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 1300
$ws.Range("D15").Value = 38000
$ws.Range("D17").Value = 493800
$ws.Range("D18").Value = 66700
$ws.Range("D20").Value = 800
$ws.Range("D21").Value = 106100
$ws.Range("D22").Value = 15600
$ws.Range("D23").Value = 51900
$ws.Range("D24").Value = 300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 51600
$ws.Range("D27").Value = 50600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -800
$ws.Range("D33").Value = 50600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 50600

# --- Balance Sheet (rows 38-77) ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 85400
$ws.Range("D42").Value = "NA"
$ws.Range("D43").Value = 37900
$ws.Range("D44").Value = "NA"
$ws.Range("D45").Value = 7800
$ws.Range("D46").Value = 131100
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 12400
$ws.Range("D49").Value = 843300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 11000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 997700
$ws.Range("D57").Value = 3800
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 44500
$ws.Range("D60").Value = 48300
$ws.Range("D61").Value = 260400
$ws.Range("D62").Value = 56500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 390100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -961700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 607600
$ws.Range("D77").Value = 0

# --- Cash Flow (rows 80-102) ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 50600
$ws.Range("D83").Value = 38600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 105500
$ws.Range("D91").Value = -7600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -7600
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -74800
$ws.Range("D101").Value = -100
$ws.Range("D102").Value = 23100

# Match column widths roughly to the neighboring (existing) columns
# now that D holds the same kind of data as E:K.
$ws.Range("D7:D102").ColumnWidth = $ws.Range("E7").ColumnWidth
